$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right
$ws.Columns("A").Insert()

# Match the styles used by the diff: header row keeps the bold/fill style
# (same as the rest of row 1, now shifted to B1), and the data row uses the
# plain left-aligned style (same as C2, not the quotePrefix style that the
# shifted-right B2 cell carries).
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("C2").Copy($ws.Range("A2"))

# Populate the new column A with the RefID header/value
$ws.Range("A1").Value = "RefID"
$ws.Range("A2").Value = "LATFLD-2"

# Select A2 to match the saved selection state
$ws.Range("A2").Select()
